# Update "想去人数" (want-to-go count) figures across sheets to match
# the latest generated output (gh-pages build at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 12522
$wsExpo.Range("F5").Value = 307
$wsExpo.Range("F8").Value = 367

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 27

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 27
$wsAll.Range("F5").Value = 12522
$wsAll.Range("F6").Value = 307
$wsAll.Range("F11").Value = 367
